# Bird feeder BOM update: add new/updated component rows (R1-R6, C1-C5,
# IC4, CN1, battery) and refresh the sheet view/selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: R1, R2, R3, R4 (qty 4) -> Yageo 10K thin film resistor ---
$ws.Cells.Item(12, 2).Value = 4
$ws.Cells.Item(12, 3).Value = "R1, R2, R3, R4"
$ws.Cells.Item(12, 4).Value = "Yageo"
$ws.Cells.Item(12, 5).Value = "RT0201FRE0710KL"
$ws.Cells.Item(12, 6).Value = "Thin Film Resistor, 10K"
$ws.Cells.Item(12, 7).Value = "Mouser"
$ws.Cells.Item(12, 8).Value = "603-RT0201FRE0710KL"

# --- Row 13: R5 -> Yageo 1K thin film resistor ---
$ws.Cells.Item(13, 3).Value = "R5"
$ws.Cells.Item(13, 4).Value = "Yageo"
$ws.Cells.Item(13, 5).Value = "RT0805FRE101KL"
$ws.Cells.Item(13, 6).Value = "Thin Film Resistor, 1K"
$ws.Cells.Item(13, 7).Value = "Mouser"
$ws.Cells.Item(13, 8).Value = "603-RT0805FRE101KL"

# --- Row 14: R6 -> Yageo 100K thin film resistor ---
$ws.Cells.Item(14, 3).Value = "R6"
$ws.Cells.Item(14, 4).Value = "Yageo"
$ws.Cells.Item(14, 5).Value = "AT0402DRE07100KL"
$ws.Cells.Item(14, 6).Value = "Thin Film Resistor, 100K"
$ws.Cells.Item(14, 7).Value = "Mouser"
$ws.Cells.Item(14, 8).Value = "603-AT0402DRE07100KL"

# --- Row 15: C1 -> KEMET 0.33uF ceramic capacitor ---
$ws.Cells.Item(15, 3).Value = "C1"
$ws.Cells.Item(15, 4).Value = "KEMET"
$ws.Cells.Item(15, 5).Value = "C0603C334K8RACAUTO"
$ws.Cells.Item(15, 6).Value = "Multilayer Ceramic Capacitor SMT, 0.33uF"
$ws.Cells.Item(15, 7).Value = "Mouser"
$ws.Cells.Item(15, 8).Value = "80-C0603C334K8RAUTO"

# --- Row 16: C2, C3, C4, C5 (qty 4) -> KEMET 0.1uF ceramic capacitor ---
$ws.Cells.Item(16, 2).Value = 4
$ws.Cells.Item(16, 3).Value = "C2, C3, C4, C5"
$ws.Cells.Item(16, 4).Value = "KEMET"
$ws.Cells.Item(16, 5).Value = "C0402C104K8PAC7411"
$ws.Cells.Item(16, 6).Value = "Multilayer Ceramic Capacitor SMT, 0.1uF"
$ws.Cells.Item(16, 7).Value = "Mouser"
$ws.Cells.Item(16, 8).Value = "80-C0402C104K8P7411"

# --- Row 17: SENSOR (Force Sensor) - unchanged content, rewritten for clarity ---
$ws.Cells.Item(17, 3).Value = "SENSOR (Force Sensor)"
$ws.Cells.Item(17, 4).Value = "Ohmite"
$ws.Cells.Item(17, 5).Value = "FSR03CE"
$ws.Cells.Item(17, 6).Value = "Force Sensing Resistor"
$ws.Cells.Item(17, 7).Value = "Mouser"
$ws.Cells.Item(17, 8).Value = "588-FSR03CE "

# --- Row 18 (new): IC4 -> Linear voltage regulator ---
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = "IC4"
$ws.Cells.Item(18, 4).Value = "Texas Instruments"
$ws.Cells.Item(18, 5).Value = "LM7805CT/NOPB"
$ws.Cells.Item(18, 6).Value = "Linear Voltage Regulator"
$ws.Cells.Item(18, 7).Value = "Mouser"
$ws.Cells.Item(18, 8).Value = "926-LM7805CT/NOPB"

# --- Row 19 (new): CN1 -> JST right-angle connector ---
$ws.Cells.Item(19, 1).Value = 12
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = "CN1"
$ws.Cells.Item(19, 4).Value = "SparkFun"
$ws.Cells.Item(19, 5).Value = "PRT-09749"
$ws.Cells.Item(19, 6).Value = "JST Right-Angle Connector TH 2-Pin"
$ws.Cells.Item(19, 7).Value = "Mouser"
$ws.Cells.Item(19, 8).Value = "474-PRT-09749"
$ws.Cells.Item(19, 5).HorizontalAlignment = -4131
$ws.Cells.Item(19, 5).VerticalAlignment = -4160

# --- Row 20 (new): 7.4V lithium battery ---
$ws.Cells.Item(20, 1).Value = 13
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = "N/A"
$ws.Cells.Item(20, 4).Value = "HHZ"
$ws.Cells.Item(20, 5).Value = "N/A"
$ws.Cells.Item(20, 6).Value = "7.4V Lithium Battery"
$ws.Cells.Item(20, 7).Value = "Amazon"
$ws.Cells.Item(20, 8).Value = "B07Q2CFRKW"

# Widen columns E and H to fit the longer MPN/VPN text now in the table.
$ws.Columns.Item(5).ColumnWidth = 20.59
$ws.Columns.Item(8).ColumnWidth = 20.92

# Scroll the view down to the table body and leave the selection where the
# editor last left off.
$ws.Range("D17").Select()

$wb.Save()
